$d = $word.ActiveDocument

# --- 1) Split the "Выполнили: ..." credits line -------------------------
# Original: "Выполнили: Шарипов Николай, Супрунов Матвей"
# Target:   "Выполнили: Шарипов Николай, " (run 1, unchanged formatting)
#         + "Изотова Валентина" (run 2, same italic Times New Roman 16pt formatting,
#           but emitted as its own run)
$oldName = "Супрунов Матвей"
$newName = "Изотова Валентина"
$fullText = $d.Content.Text
$idx = $fullText.IndexOf($oldName)

if ($idx -ge 0) {
    $target = $d.Range($idx, $idx + $oldName.Length)
    $target.Text = ""

    $insertionPoint = $d.Range($idx, $idx)
    $insertionPoint.InsertAfter($newName)

    $newRange = $d.Range($idx, $idx + $newName.Length)
    # Re-assert the (already identical) character formatting on the freshly
    # typed text so Word keeps it as a discrete run instead of silently
    # re-merging it into the preceding run.
    $newRange.Font.Bold = $true
    $newRange.Font.Bold = $false
    $newRange.Font.Italic = $true
    $newRange.Font.Size = 16
    $newRange.Font.NameBi = "Times New Roman"
}

# --- 2) Mark every inline picture's run as NoProofing (w:noProof) -------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}

Write-Output "done"
